$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 2000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2000
$ws.Range("N13").Value = -2338
$ws.Range("H34").Value = 994.8
$ws.Range("H36").Value = 994.8
$ws.Range("H61").Value = 28727.5
$ws.Range("I61").Value = 28727.5
$ws.Range("K61").Value = 86182.5
$ws.Range("M61").Value = -86010.5
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480
$ws.Range("H113").Value = 6970.0713
$ws.Range("I113").Value = 6599.8
$ws.Range("J113").Value = 7175.778
$ws.Range("K113").Value = 6599.8
$ws.Range("L113").Value = 7175.778
$ws.Range("M113").Value = -3345.8
$ws.Range("N113").Value = -13683.778
$ws.Range("H137").Value = 3638.05
$ws.Range("I137").Value = 3224.3333
$ws.Range("J137").Value = 3976.5454
$ws.Range("K137").Value = 9672.999899999999
$ws.Range("L137").Value = 11929.6362
$ws.Range("M137").Value = -7122.999899999999
$ws.Range("N137").Value = -17029.6362
$ws.Range("M13").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2277.6
$ws.Range("I45").Value = 1796
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 1796
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -1419
$ws.Range("N45").Value = -3754
$ws.Range("H61").Value = 3278.8
$ws.Range("J61").Value = 4744
$ws.Range("L61").Value = 4744
$ws.Range("N61").Value = -5168
$ws.Range("H110").Value = 2241.6667
$ws.Range("I110").Value = 2146.875
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 2146.875
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = -101.875
$ws.Range("N110").Value = -7090
$ws.Range("H136").Value = 3278.8
$ws.Range("J136").Value = 4744
$ws.Range("L136").Value = 14232
$ws.Range("N136").Value = -19332

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2269.818
$ws.Range("I58").Value = 2219.3333
$ws.Range("K58").Value = 2219.3333
$ws.Range("M58").Value = -2016.3333
$ws.Range("H105").Value = 632.8333
$ws.Range("I105").Value = 626.8182
$ws.Range("K105").Value = 626.8182
$ws.Range("M105").Value = 1120.1818
$ws.Range("H136").Value = 2269.818
$ws.Range("I136").Value = 2219.3333
$ws.Range("K136").Value = 6657.999899999999
$ws.Range("M136").Value = -4107.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 278
$ws.Range("I7").Value = 237
$ws.Range("K7").Value = 711
$ws.Range("M7").Value = -599
$ws.Range("H25").Value = 1665.25
$ws.Range("I25").Value = 886.6
$ws.Range("J25").Value = 2963
$ws.Range("K25").Value = 2659.8
$ws.Range("L25").Value = 8889
$ws.Range("M25").Value = -2490.8
$ws.Range("N25").Value = -9227
$ws.Range("H30").Value = 1665.25
$ws.Range("I30").Value = 886.6
$ws.Range("J30").Value = 2963
$ws.Range("K30").Value = 2659.8
$ws.Range("L30").Value = 8889
$ws.Range("M30").Value = -2557.8
$ws.Range("N30").Value = -9093
$ws.Range("H46").Value = 322
$ws.Range("I46").Value = 444
$ws.Range("K46").Value = 1332
$ws.Range("M46").Value = -1241
$ws.Range("H80").Value = 10995.5
$ws.Range("J80").Value = 10995.5
$ws.Range("L80").Value = 32986.5
$ws.Range("N80").Value = -34858.5
$ws.Range("H83").Value = 10995.5
$ws.Range("J83").Value = 10995.5
$ws.Range("L83").Value = 98959.5
$ws.Range("N83").Value = -108319.5
$ws.Range("H92").Value = 1067.5
$ws.Range("J92").Value = 1123.3334
$ws.Range("L92").Value = 3370.0002
$ws.Range("N92").Value = -5866.0002
$ws.Range("H131").Value = 1374.7
$ws.Range("I131").Value = 944.75
$ws.Range("J131").Value = 1661.3334
$ws.Range("K131").Value = 2834.25
$ws.Range("L131").Value = 4984.0002
$ws.Range("M131").Value = 2205.75
$ws.Range("N131").Value = -15064.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3731.6316
$ws.Range("I43").Value = 2456.2727
$ws.Range("J43").Value = 5485.25
$ws.Range("K43").Value = 2456.2727
$ws.Range("L43").Value = 5485.25
$ws.Range("M43").Value = -2305.2727
$ws.Range("N43").Value = -5787.25
$ws.Range("H57").Value = 19765.25
$ws.Range("J57").Value = 19765.25
$ws.Range("L57").Value = 19765.25
$ws.Range("N57").Value = -21405.25
$ws.Range("H102").Value = 2457.1875
$ws.Range("I102").Value = 2457.1875
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2457.1875
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -835.1875
$ws.Range("H113").Value = 1555.5
$ws.Range("I113").Value = 1111
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1111
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1059
$ws.Range("N113").Value = -6340
$ws.Range("N102").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3131.5
$ws.Range("I122").Value = 3113
$ws.Range("J122").Value = 3150
$ws.Range("K122").Value = 9339
$ws.Range("L122").Value = 9450
$ws.Range("M122").Value = -6889
$ws.Range("N122").Value = -14350
$ws.Range("H136").Value = 3763.1667
$ws.Range("I136").Value = 2815.8572
$ws.Range("J136").Value = 5089.4
$ws.Range("K136").Value = 8447.571599999999
$ws.Range("L136").Value = 15268.2
$ws.Range("M136").Value = -5897.571599999999
$ws.Range("N136").Value = -20368.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 993
$ws.Range("I81").Value = 791.6
$ws.Range("K81").Value = 1583.2
$ws.Range("M81").Value = -522.2
$ws.Range("H84").Value = 993
$ws.Range("I84").Value = 791.6
$ws.Range("K84").Value = 7916
$ws.Range("M84").Value = -2612
$ws.Range("H107").Value = 1039.1666
$ws.Range("I107").Value = 733.6667
$ws.Range("J107").Value = 2566.6667
$ws.Range("K107").Value = 2201.0001
$ws.Range("L107").Value = 7700.000100000001
$ws.Range("M107").Value = -281.0001000000002
$ws.Range("N107").Value = -11540.0001
$ws.Range("H126").Value = 22043.777
$ws.Range("I126").Value = 23431.105
$ws.Range("J126").Value = 18748.875
$ws.Range("K126").Value = 70293.315
$ws.Range("L126").Value = 56246.625
$ws.Range("M126").Value = -67823.315
$ws.Range("N126").Value = -61186.625
$ws.Range("H136").Value = 8664.833000000001
$ws.Range("I136").Value = 8542.429
$ws.Range("J136").Value = 8836.200000000001
$ws.Range("K136").Value = 25627.287
$ws.Range("L136").Value = 26508.6
$ws.Range("M136").Value = -23077.287
$ws.Range("N136").Value = -31608.6
